$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data submitted via the backend commit functionality
$newRows = @(
    @{ Row = 16; DateTime = "2025-05-05T14:57:30.993Z"; Ip = "149.132.26.74"; C = 1; D = 1; E = 0; F = 1; G = 0; Prediction = "uncertain case" },
    @{ Row = 17; DateTime = "2025-05-07T16:19:03.014Z"; Ip = "93.66.5.34";    C = 0; D = 0; E = 0; F = 0; G = 0; Prediction = "Likely NOT Malignant" },
    @{ Row = 18; DateTime = "2025-05-07T16:22:36.511Z"; Ip = "93.66.5.34";    C = 1; D = 1; E = 0; F = 0; G = 0; Prediction = "Likely Malignant" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.DateTime
    $ws.Cells.Item($row, 2).Value = $r.Ip
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.Prediction
}
